# Refresh of the live cryptocurrency price/volume snapshot on Sheet1
# (scheduled scrape update, mirrors the "Updated cryptos list ... with GitHub
# Actions" commit). Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Rows 15/16 and 38/39 additionally swap which coin occupies which row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.586.47'
$ws.Range('E2').Value = '  +1.28%  '

# Row 3
$ws.Range('D3').Value = '2.624.77'
$ws.Range('E3').Value = '  +0.12%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').Value = '''595.01'
$ws.Range('E5').Value = '  -0.37%  '

# Row 6
$ws.Range('D6').Value = '''152.67'
$ws.Range('E6').Value = '  +1.16%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').Value = '''0.590'
$ws.Range('E8').Value = '  -0.16%  '

# Row 9
$ws.Range('D9').Value = '''0.115'
$ws.Range('E9').Value = '  +4.90%  '

# Row 10
$ws.Range('D10').Value = '''0.397'
$ws.Range('E10').Value = '  +3.27%  '

# Row 11
$ws.Range('D11').Value = '''5.81'
$ws.Range('E11').Value = '  +1.96%  '

# Row 12
$ws.Range('D12').Value = '''0.153'
$ws.Range('E12').Value = '  +1.27%  '

# Row 13
$ws.Range('D13').Value = '''28.69'
$ws.Range('E13').Value = '  +3.09%  '

# Row 14
$ws.Range('D14').Value = '3.094.30'
$ws.Range('E14').Value = '  +0.13%  '

# Row 15
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '''0.0000172'
$ws.Range('E15').Value = '  +12.15%  '

# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '64.420.17'
$ws.Range('E16').Value = '  +1.28%  '

# Row 17
$ws.Range('D17').Value = '2.581.79'
$ws.Range('E17').Value = '  -1.50%  '

# Row 18
$ws.Range('D18').Value = '''12.29'
$ws.Range('E18').Value = '  -0.39%  '

# Row 19
$ws.Range('D19').Value = '''4.79'
$ws.Range('E19').Value = '  +1.85%  '

# Row 20
$ws.Range('D20').Value = '''351.53'
$ws.Range('E20').Value = '  +1.07%  '

# Row 21
$ws.Range('D21').Value = '''7.18'

# Row 22
$ws.Range('E22').Value = '  +0.13%  '

# Row 23
$ws.Range('D23').Value = '''67.69'

# Row 24
$ws.Range('D24').Value = '''1.70'
$ws.Range('E24').Value = '  -1.03%  '

# Row 25
$ws.Range('D25').Value = '''9.32'
$ws.Range('E25').Value = '  +0.15%  '

# Row 26
$ws.Range('D26').Value = '''1.65'
$ws.Range('E26').Value = '  -1.39%  '

# Row 27
$ws.Range('D27').Value = '''8.26'
$ws.Range('E27').Value = '  +1.13%  '

# Row 28
$ws.Range('D28').Value = '''0.165'
$ws.Range('E28').Value = '  +1.66%  '

# Row 29
$ws.Range('D29').Value = '''548.21'
$ws.Range('E29').Value = '  -1.35%  '

# Row 30
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.02%  '

# Row 31
$ws.Range('D31').Value = '0.0₃0917'
$ws.Range('E31').Value = '  +7.78%  '

# Row 32
$ws.Range('E32').Value = '  +0.96%  '

# Row 33
$ws.Range('D33').Value = '''1.82'
$ws.Range('E33').Value = '  +2.66%  '

# Row 34
$ws.Range('D34').Value = '''5.72'
$ws.Range('E34').Value = '  +8.31%  '

# Row 35
$ws.Range('D35').Value = '''6.23'
$ws.Range('E35').Value = '  +0.36%  '

# Row 36
$ws.Range('D36').Value = '''0.423'
$ws.Range('E36').Value = '  +2.03%  '

# Row 37
$ws.Range('D37').Value = '''164.06'
$ws.Range('E37').Value = '  -2.38%  '

# Row 38
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '''2.01'
$ws.Range('E38').Value = '  +3.62%  '

# Row 39
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = '''20.12'
$ws.Range('E39').Value = '  +2.99%  '

# Row 40
$ws.Range('D40').Value = '''0.998'
$ws.Range('E40').Value = '  -0.20%  '

# Row 41
$ws.Range('E41').Value = '  -0.03%  '

# Row 42
$ws.Range('D42').Value = '''168.91'
$ws.Range('E42').Value = '  +1.36%  '

# Row 43
$ws.Range('D43').Value = '''41.92'
$ws.Range('E43').Value = '  +5.47%  '

# Row 44
$ws.Range('D44').Value = '''4.09'
$ws.Range('E44').Value = '  +4.03%  '

# Row 45
$ws.Range('D45').Value = '''23.34'
$ws.Range('E45').Value = '  +7.59%  '

# Row 46
$ws.Range('D46').Value = '''0.0599'
$ws.Range('E46').Value = '  +1.42%  '

# Row 47
$ws.Range('D47').Value = '''2.24'
$ws.Range('E47').Value = '  +11.04%  '

# Row 48
$ws.Range('D48').Value = '''0.641'
$ws.Range('E48').Value = '  +1.60%  '

# Row 49
$ws.Range('D49').Value = '''0.0252'
$ws.Range('E49').Value = '  -0.04%  '

# Row 50
$ws.Range('D50').Value = '''0.0983'
$ws.Range('E50').Value = '  +1.69%  '

# Row 51
$ws.Range('D51').Value = '''19.33'
$ws.Range('E51').Value = '  -0.13%  '
